# Sales Order Sales added and Receipt fields added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")

# Copy formatting from the row above (row 26) into the new row (row 27)
# so the new row visually matches the rest of the table (8pt font, etc).
$ws.Range("A26:E26").Copy($ws.Range("A27:E27"))

# Populate the new "sales_order_sales" table row
$ws.Range("A27").Value = "sales_order_sales"
$ws.Range("B27").Value = "so,product,quantity,transaction,sale_quantity"
$ws.Range("C27").Value = "USER"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 10

# Grow the worksheet table ("Table2") so the new row is included,
# which also keeps the AutoFilter range in sync.
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A1:E27"))
